$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1. Title heading ("Play Megajackpots Ocean Belles Free - Pin-Up Inspired Slot
#    Game" -> "Play Megajackpots Ocean Belles for Free"). This single Find/Replace
#    call (scoped to the whole document) also fixes the matching bold call-to-
#    action paragraph further down that repeats the exact same sentence.
$d.Content.Find.Execute("Play Megajackpots Ocean Belles Free - Pin-Up Inspired Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play Megajackpots Ocean Belles for Free", 2)

# 2. "What we like" / "What we don't like" bullet lists - each bullet
#    paragraph starts with an empty run (<w:r/>) ahead of the text run; use
#    InsertXML with the reconstructed paragraph so that leading empty run is
#    preserved exactly as it is in the source (plain Find/Replace merges it
#    away because the replacement text run carries no <w:rPr>).
function Set-BulletParagraph($searchText, $newText) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $p = $rng.Paragraphs(1)
    $xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:spacing w:line=`"240`" w:lineRule=`"auto`"/><w:ind w:left=`"720`"/></w:pPr><w:r/><w:r><w:t>$newText</w:t></w:r></w:p>"
    $p.Range.InsertXML($xml) | Out-Null
}

Set-BulletParagraph "Beautiful visuals with a pin-up theme" "Unique and charming pin-up theme"
Set-BulletParagraph "Game structure is unique and engaging" "Beautiful visuals set on a luxurious ship"
Set-BulletParagraph "Wild symbols can expand to cover entire reels" "Expanding wild symbols for bigger wins"
Set-BulletParagraph "Bonus feature initiates free spins" "Engaging bonus rounds with free spins"

# 3. "What we don't like" bullets
Set-BulletParagraph "Theme may not appeal to everyone" "Limited availability of slot games with pin-up theme"
Set-BulletParagraph "75 pay lines may be overwhelming for some players" "Some players may prefer different game structures"

# 4. Closing italic meta description
$d.Content.Find.Execute("Read our review of Megajackpots Ocean Belles, a pin-up inspired slot game with bonus features. Play for free and experience the charm of 1950s pin-ups.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Megajackpots Ocean Belles slot game and play for free. Enjoy the unique pin-up theme and engaging bonus features.", 2)
